# Normalise punctuation in a handful of reference/description cells on the
# "DBD" sheet (full-width colons / stray periods after the leading code
# -> plain half-width colons), and reflow the "00A/201" book-code note onto
# two lines. Finishes with DBD as the active sheet/selection, matching the
# author's last on-screen position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("G24").Value = "CdAcCode 會計科子細目設定檔`nIC1:短擔息`nIC2:中擔息`nIC3:長擔息`nIC4:三十年房貸息"

$ws.Range("G27").Value = "以放款主檔的下次應繳日~本營業日計算`n0:一個月以下`n1:一～三個月`n2:三～六個月`n3:六個月以上`n"

$ws.Range("G28").Value = "000:全公司"

$ws.Range("G29").Value = "00A:傳統帳冊`n201:利變年金帳冊"

$null = $ws.Activate()
$null = $ws.Range("H28").Select()
